$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Date Published changes from 2020-07-21 (44033) to 2020-07-22 (44034)
$ws.Range("B3").Value = 44034

# Row 4 (New York -- New York): the scrape for this run came back empty for
# most columns and flipped "Pct Includes Hispanic Black" (J4) to FALSE.
# B4 also loses its date formatting, so clear format+contents together.
$ws.Range("B4").Clear()
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""

# Row 39 (Delaware): updated error message in Status code column
$ws.Range("O39").Value = "An error occurred. ... HTTPError('504 Server Error: Gateway Time-out for url: https://myhealthycommunity.dhss.delaware.gov/locations/state/')"
